$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "All Trades" sheet: append trades #82 (momentum DOWN) and #83
# (HighProbConvergence UP) as new rows 83 and 84.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Seed the new rows from the last existing row of the same shape (row 82) so
# that column types/formatting (text dates, empty-string placeholder cells
# for Exit Price / Exit Reason, etc.) carry over exactly, then overwrite just
# the cells that actually differ for each new trade.
$wsAll.Rows.Item(82).Copy()
$wsAll.Rows.Item(83).PasteSpecial(-4163)
$wsAll.Cells.Item(83, 1).Value = 82
$wsAll.Cells.Item(83, 3).Value = "00:12:45"
$wsAll.Cells.Item(83, 4).Value = "momentum"
$wsAll.Cells.Item(83, 5).Value = "DOWN"
$wsAll.Cells.Item(83, 6).Value = 0.66
$wsAll.Cells.Item(83, 16).Value = 0.9
$wsAll.Cells.Item(83, 17).Value = "Downward momentum: -1.980% over 10 samples"

$wsAll.Rows.Item(82).Copy()
$wsAll.Rows.Item(84).PasteSpecial(-4163)
$wsAll.Cells.Item(84, 1).Value = 83
$wsAll.Cells.Item(84, 3).Value = "00:12:45"
$wsAll.Cells.Item(84, 4).Value = "HighProbConvergence"
$wsAll.Cells.Item(84, 5).Value = "UP"
$wsAll.Cells.Item(84, 6).Value = 0.33
$wsAll.Cells.Item(84, 16).Value = 0.95
$wsAll.Cells.Item(84, 17).Value = "Mean reversion UP: price 1.59% below mean (z=-2.00)"

# ---------------------------------------------------------------------------
# "momentum" strategy sheet: append the same trade #82 as a new row 14.
# ---------------------------------------------------------------------------
$wsMom = $wb.Worksheets.Item("momentum")

$wsMom.Rows.Item(13).Copy()
$wsMom.Rows.Item(14).PasteSpecial(-4163)
$wsMom.Cells.Item(14, 1).Value = 82
$wsMom.Cells.Item(14, 3).Value = "00:12:45"
$wsMom.Cells.Item(14, 6).Value = 0.66
$wsMom.Cells.Item(14, 12).Value = 0
$wsMom.Cells.Item(14, 13).Value = 0
$wsMom.Cells.Item(14, 14).Value = 0.9
$wsMom.Cells.Item(14, 15).Value = "Downward momentum: -1.980% over 10 samples"
$wsMom.Cells.Item(14, 17).Value = 0

# ---------------------------------------------------------------------------
# "HighProbConvergence" strategy sheet: append trade #83 as a new row 8.
# ---------------------------------------------------------------------------
$wsHpc = $wb.Worksheets.Item("HighProbConvergence")

$wsHpc.Rows.Item(7).Copy()
$wsHpc.Rows.Item(8).PasteSpecial(-4163)
$wsHpc.Cells.Item(8, 1).Value = 83
$wsHpc.Cells.Item(8, 3).Value = "00:12:45"
$wsHpc.Cells.Item(8, 6).Value = 0.33
$wsHpc.Cells.Item(8, 12).Value = 0
$wsHpc.Cells.Item(8, 13).Value = 0
$wsHpc.Cells.Item(8, 14).Value = 0.95
$wsHpc.Cells.Item(8, 15).Value = "Mean reversion UP: price 1.59% below mean (z=-2.00)"
$wsHpc.Cells.Item(8, 17).Value = 0
